$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H44").Value = 33025
$ws.Range("J44").Value = 33025
$ws.Range("L44").Value = 33025
$ws.Range("N44").Value = -33949
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()
$ws.Range("H88").Value = 2874.25
$ws.Range("I88").Value = 2874.5
$ws.Range("J88").Value = 2874
$ws.Range("K88").Value = 2874.5
$ws.Range("L88").Value = 2874
$ws.Range("M88").Value = -2468.5
$ws.Range("N88").Value = -3686
$ws.Range("H91").Value = 2874.25
$ws.Range("I91").Value = 2874.5
$ws.Range("J91").Value = 2874
$ws.Range("K91").Value = 2874.5
$ws.Range("L91").Value = 2874
$ws.Range("M91").Value = -1470.5
$ws.Range("N91").Value = -5682
$ws.Range("H105").Value = 16791.143
$ws.Range("J105").Value = 16791.143
$ws.Range("L105").Value = 16791.143
$ws.Range("N105").Value = -23779.143

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H101").Value = 61700.5
$ws.Range("J101").Value = 61700.5
$ws.Range("L101").Value = 61700.5
$ws.Range("N101").Value = -68190.5
$ws.Range("H102").Value = 15627546
$ws.Range("I102").Value = 17858782
$ws.Range("K102").Value = 17858782
$ws.Range("M102").Value = -17857160

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 83334490
$ws.Range("J107").Value = 500
$ws.Range("L107").Value = 500
$ws.Range("N107").Value = -4340
$ws.Range("H132").Value = 125116.664
$ws.Range("J132").Value = 125116.664
$ws.Range("L132").Value = 125116.664
$ws.Range("N132").Value = -135236.664

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H33").Value = 1149.5
$ws.Range("I33").Value = 1149.5
$ws.Range("K33").Value = 1149.5
$ws.Range("M33").Value = -770.5
$ws.Range("H35").Value = 1002.8571
$ws.Range("I35").Value = 1002.8571
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 1002.8571
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -708.8570999999999
$ws.Range("N35").ClearContents()
$ws.Range("H42").Value = 14500
$ws.Range("I42").Value = 14500
$ws.Range("K42").Value = 14500
$ws.Range("M42").Value = -13907
$ws.Range("H44").Value = 2500
$ws.Range("J44").Value = 2500
$ws.Range("L44").Value = 2500
$ws.Range("N44").Value = -3384
$ws.Range("H55").Value = 14000
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()
$ws.Range("H68").Value = 88382.5
$ws.Range("J68").Value = 88382.5
$ws.Range("L68").Value = 88382.5
$ws.Range("N68").Value = -89880.5
$ws.Range("H71").Value = 88382.5
$ws.Range("J71").Value = 88382.5
$ws.Range("L71").Value = 265147.5
$ws.Range("N71").Value = -272635.5
$ws.Range("H93").Value = 8063.857
$ws.Range("I93").Value = 4333.1665
$ws.Range("K93").Value = 4333.1665
$ws.Range("M93").Value = -2461.1665
$ws.Range("H105").Value = 2971.6667
$ws.Range("I105").Value = 1346.8
$ws.Range("K105").Value = 1346.8
$ws.Range("M105").Value = 400.2
$ws.Range("H107").Value = 1625.8462
$ws.Range("I107").Value = 607.875
$ws.Range("J107").Value = 3254.6
$ws.Range("K107").Value = 607.875
$ws.Range("L107").Value = 3254.6
$ws.Range("M107").Value = 1312.125
$ws.Range("N107").Value = -7094.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 241.33333
$ws.Range("J97").Value = 152
$ws.Range("L97").Value = 456
$ws.Range("N97").Value = -1448
$ws.Range("H119").Value = 6833.3335
$ws.Range("I119").Value = 6833.3335
$ws.Range("K119").Value = 20500.0005
$ws.Range("M119").Value = -15662.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 60470.668
$ws.Range("J57").Value = 70706
$ws.Range("L57").Value = 70706
$ws.Range("N57").Value = -72346
$ws.Range("H80").Value = 2099.8
$ws.Range("J80").Value = 2499.5
$ws.Range("L80").Value = 2499.5
$ws.Range("N80").Value = -4495.5
$ws.Range("H83").Value = 2099.8
$ws.Range("J83").Value = 2499.5
$ws.Range("L83").Value = 12497.5
$ws.Range("N83").Value = -22481.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1333.1666
$ws.Range("I22").Value = 1000
$ws.Range("J22").Value = 1499.75
$ws.Range("K22").Value = 1000
$ws.Range("L22").Value = 1499.75
$ws.Range("M22").Value = -705
$ws.Range("N22").Value = -2089.75
$ws.Range("H27").Value = 1333.1666
$ws.Range("I27").Value = 1000
$ws.Range("J27").Value = 1499.75
$ws.Range("K27").Value = 1000
$ws.Range("L27").Value = 1499.75
$ws.Range("M27").Value = -893
$ws.Range("N27").Value = -1713.75
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("H38").Value = 8500
$ws.Range("J38").Value = 8500
$ws.Range("L38").Value = 8500
$ws.Range("N38").Value = -9320
$ws.Range("H61").Value = 333338000
$ws.Range("I61").Value = 333338000
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 333338000
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -333337798
$ws.Range("N61").ClearContents()
$ws.Range("H68").Value = 5798
$ws.Range("J68").Value = 5249.5
$ws.Range("L68").Value = 5249.5
$ws.Range("N68").Value = -6747.5
$ws.Range("H71").Value = 5798
$ws.Range("J71").Value = 5249.5
$ws.Range("L71").Value = 26247.5
$ws.Range("N71").Value = -33735.5
$ws.Range("H82").Value = 694.55554
$ws.Range("I82").Value = 719.2
$ws.Range("J82").Value = 663.75
$ws.Range("K82").Value = 719.2
$ws.Range("L82").Value = 663.75
$ws.Range("M82").Value = -358.2
$ws.Range("N82").Value = -1385.75
$ws.Range("H85").Value = 694.55554
$ws.Range("I85").Value = 719.2
$ws.Range("J85").Value = 663.75
$ws.Range("K85").Value = 719.2
$ws.Range("L85").Value = 663.75
$ws.Range("M85").Value = 528.8
$ws.Range("N85").Value = -3159.75
$ws.Range("H113").Value = 333338000
$ws.Range("I113").Value = 333338000
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 333338000
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -333335830
$ws.Range("N113").ClearContents()
$ws.Range("H136").Value = 3549.8333
$ws.Range("I136").Value = 3399.8
$ws.Range("J136").Value = 4300
$ws.Range("K136").Value = 10199.4
$ws.Range("L136").Value = 12900
$ws.Range("M136").Value = -7649.400000000001
$ws.Range("N136").Value = -18000

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H30").Value = 4850
$ws.Range("J30").Value = 4850
$ws.Range("L30").Value = 4850
$ws.Range("N30").Value = -5064
$ws.Range("H34").Value = 25513
$ws.Range("I34").Value = 31026
$ws.Range("J34").Value = 20000
$ws.Range("K34").Value = 31026
$ws.Range("L34").Value = 20000
$ws.Range("M34").Value = -30823
$ws.Range("N34").Value = -20406
$ws.Range("H132").Value = 2997
$ws.Range("I132").Value = 2994
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 8982
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -6452
$ws.Range("N132").Value = -14060
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H136").Value = 2679.8572
$ws.Range("I136").Value = 1762.6
$ws.Range("K136").Value = 5287.799999999999
$ws.Range("M136").Value = -2737.799999999999
